$wb = $excel.ActiveWorkbook

# Sheet: AR
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.02999222118014286
$ws.Range("B3").Value = 0.7325608976808113
$ws.Range("B4").Value = 0.1408722516562999

# Sheet: SETAR
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B2").Value = -0.1974865599909555
$ws.Range("B3").Value = 0.4714135928410155
$ws.Range("B4").Value = 0.08615543052135079
$ws.Range("B5").Value = 0.2212855327134672
$ws.Range("B6").Value = 0.5199371156345707
$ws.Range("B7").Value = 0.11074069437745

# Sheet: GARCH
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.001719533569652933
$ws.Range("B3").Value = 0.123955123772953
$ws.Range("B4").Value = 0.1171079610609286
$ws.Range("B5").Value = 0.09761581642944596

# Sheet: TARCH
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = 0.005193054085146191
$ws.Range("B3").Value = 0.126913691658619
$ws.Range("B4").Value = 0.1624302673115644
$ws.Range("B5").Value = -0.09557983522698058
$ws.Range("B6").Value = 0.08057397974730415

# Sheet: AR-TARCH
$ws = $wb.Worksheets.Item("AR-TARCH")
$ws.Range("B2").Value = 0.02991105400350757
$ws.Range("B3").Value = 0.7421410095620198
$ws.Range("B4").Value = 0.122494724306115
$ws.Range("B5").Value = 0.1600111935395574
$ws.Range("B6").Value = -0.08746286348230323
$ws.Range("B7").Value = 2.06478559557175 * [Math]::Pow(10, -13)
